$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: new (empty) cell N2, formatted like M2 ---
$ws.Range("M2").Copy()
$ws.Range("N2").PasteSpecial(-4122)

# --- Row 3: new header cell N3 = 2021, formatted like M3 ---
$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial(-4122)
$ws.Range("N3").Value = 2021

# --- Row 4: D4:L4 become bold (matching M4's existing bold style),
#            new cell N4 formatted like M4 ---
$ws.Range("D4:L4").Font.Bold = $true
$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("N4").Value = 95.134712433469176

# --- Row 5: new cell N5 formatted like M5 ---
$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial(-4122)
$ws.Range("N5").Value = 99.705541665880986

# --- Row 6: new cell N6 formatted like M6 ---
$ws.Range("M6").Copy()
$ws.Range("N6").PasteSpecial(-4122)
$ws.Range("N6").Value = 92.425193326577897

# --- Row 7: new cell N7 formatted like D7 ---
$ws.Range("D7").Copy()
$ws.Range("N7").PasteSpecial(-4122)
$ws.Range("N7").Value = 88.209991167538519

# --- Row 8: new cell N8 formatted like D8 ---
$ws.Range("D8").Copy()
$ws.Range("N8").PasteSpecial(-4122)
$ws.Range("N8").Value = 92.225038985690773

# --- Row 9: new cell N9 formatted like D9 ---
$ws.Range("D9").Copy()
$ws.Range("N9").PasteSpecial(-4122)
$ws.Range("N9").Value = 96.801032063987265

# --- Row 10: new cell N10 formatted like D10 ---
$ws.Range("D10").Copy()
$ws.Range("N10").PasteSpecial(-4122)
$ws.Range("N10").Value = 97.660491031729507

# --- Row 11: new cell N11 formatted like D11 ---
$ws.Range("D11").Copy()
$ws.Range("N11").PasteSpecial(-4122)
$ws.Range("N11").Value = 90.23262877800066

# --- Row 12: new cell N12 formatted like D12 ---
$ws.Range("D12").Copy()
$ws.Range("N12").PasteSpecial(-4122)
$ws.Range("N12").Value = 99.653994395099105

# --- Row 13: new cell N13 formatted like D13 ---
$ws.Range("D13").Copy()
$ws.Range("N13").PasteSpecial(-4122)
$ws.Range("N13").Value = 100

# --- Row 14: new cell N14 formatted like D14 ---
$ws.Range("D14").Copy()
$ws.Range("N14").PasteSpecial(-4122)
$ws.Range("N14").Value = 100

# --- Row 15: new cell N15 formatted like M15 ---
$ws.Range("M15").Copy()
$ws.Range("N15").PasteSpecial(-4122)
$ws.Range("N15").Value = 100

# --- Update print quality / vertical DPI to match new pageSetup ---
$ws.PageSetup.VerticalDpi = 300

# --- Move the active selection to N2, matching the recorded cursor position ---
$ws.Range("N2").Select() | Out-Null
